$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the E2 coefficient (Ki) value - calibrated voltage outputs
$ws.Range("E2").Value = 40000

# Force recalculation of dependent formulas (B1, B2)
$excel.Calculate()

# Move selection/active cell to E2 to match the saved view state
$ws.Range("E2").Select()
